# Update NATMI edge-weight table with recalculated receptor TPM-derived values.
# The "Target cluster" = ECs rows (2 and 5) had their receptor-expressing-cell
# counts recalculated from new TPM input, and every value derived from them
# (receptor avg/total expression, specificity columns, and the dependent edge
# weight / edge specificity columns for all six data rows) is refreshed below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending=FAPs, Target=ECs) ---
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07409833333333334
$ws.Range("N2").Value = 0.222295
$ws.Range("O2").Value = 0.02892330792240624
$ws.Range("P2").Value = 0.02892330792240624
$ws.Range("Q2").Value = 0.02929378810555556
$ws.Range("R2").Value = 0.26364409295
$ws.Range("S2").Value = 0.01600731151777644
$ws.Range("T2").Value = 0.01600731151777644

# --- Row 3 (Sending=FAPs, Target=FAPs) ---
$ws.Range("O3").Value = 0.1644833827109413
$ws.Range("P3").Value = 0.1644833827109413
$ws.Range("S3").Value = 0.09103166047310936
$ws.Range("T3").Value = 0.09103166047310939

# --- Row 4 (Sending=FAPs, Target=MuSCs) ---
$ws.Range("O4").Value = 0.8065933093666526
$ws.Range("P4").Value = 0.8065933093666525
$ws.Range("S4").Value = 0.446400889062349
$ws.Range("T4").Value = 0.446400889062349

# --- Row 5 (Sending=MuSCs, Target=ECs) ---
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07409833333333334
$ws.Range("N5").Value = 0.222295
$ws.Range("O5").Value = 0.02892330792240624
$ws.Range("P5").Value = 0.02892330792240624
$ws.Range("Q5").Value = 0.02363660265055556
$ws.Range("R5").Value = 0.212729423855
$ws.Range("S5").Value = 0.01291599640462981
$ws.Range("T5").Value = 0.01291599640462981

# --- Row 6 (Sending=MuSCs, Target=FAPs) ---
$ws.Range("O6").Value = 0.1644833827109413
$ws.Range("P6").Value = 0.1644833827109413
$ws.Range("S6").Value = 0.07345172223783189
$ws.Range("T6").Value = 0.07345172223783189

# --- Row 7 (Sending=MuSCs, Target=MuSCs) ---
$ws.Range("O7").Value = 0.8065933093666526
$ws.Range("P7").Value = 0.8065933093666525
$ws.Range("S7").Value = 0.3601924203043036
$ws.Range("T7").Value = 0.3601924203043035
